$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D9"  = -8.158999999999997
    "D18" = -8.924399999999993
    "D20" = -7.244200000000002
    "D27" = -7.825799999999998
    "D35" = -7.328800000000006
    "D69" = -7.482599999999996
    "D76" = -7.685899999999997
    "D78" = -7.557400000000001
    "D82" = -8.422600000000008
    "D83" = -8.817399999999999
    "D93" = -7.008199999999992
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
